$wb = $excel.ActiveWorkbook

# --- Sheet2: add header row + one data row ---
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("A1").Value = "Height"
$ws2.Range("B1").Value = "Weight"
$ws2.Range("C1").Value = "Shoe_size"

$ws2.Range("A2").Value = 190
$ws2.Range("B2").Value = 70
$ws2.Range("C2").Value = 43

# --- Sheet1: update selection, no longer the active tab ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$null = $ws1.Range("A1:C1").Select()

# --- Sheet2: becomes the active tab, with its own selection ---
$null = $ws2.Activate()
$null = $ws2.Range("B8").Select()
